# Weekly fruit/vegetable price update: insert 3 new price rows
# (2021-10-20 readings for "Cuatro cascos verde", "Zafiro rojo" and
# "Zafiro verde") into the existing "Pimiento" table, pushing the
# pre-existing rows 124-188 down to 127-191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 124, shifting the rest of the
# table (and its formatting) downward - mirrors Excel's own
# Insert Sheet Rows behaviour (xlShiftDown = -4121).
$ws.Range("A124:R126").EntireRow.Insert(-4121)

# Row 124: Cuatro cascos verde, 2021-10-20
$ws.Cells.Item(124, 1).Value = 11
$ws.Cells.Item(124, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(124, 3).Value = "Bíobío"
$ws.Cells.Item(124, 4).Value = 44489
$ws.Cells.Item(124, 5).Value = 8
$ws.Cells.Item(124, 6).Value = 100112002
$ws.Cells.Item(124, 7).Value = "Pimiento"
$ws.Cells.Item(124, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 100
$ws.Cells.Item(124, 11).Value = 35000
$ws.Cells.Item(124, 12).Value = 36000
$ws.Cells.Item(124, 13).Value = 35500
$ws.Cells.Item(124, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(124, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(124, 16).Value = 1972
$ws.Cells.Item(124, 17).Value = 18
$ws.Cells.Item(124, 18).Value = "Hortaliza"

# Row 125: Zafiro rojo, 2021-10-20
$ws.Cells.Item(125, 1).Value = 11
$ws.Cells.Item(125, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(125, 3).Value = "Bíobío"
$ws.Cells.Item(125, 4).Value = 44489
$ws.Cells.Item(125, 5).Value = 8
$ws.Cells.Item(125, 6).Value = 100112002
$ws.Cells.Item(125, 7).Value = "Pimiento"
$ws.Cells.Item(125, 8).Value = "Zafiro rojo"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 100
$ws.Cells.Item(125, 11).Value = 40000
$ws.Cells.Item(125, 12).Value = 42000
$ws.Cells.Item(125, 13).Value = 41000
$ws.Cells.Item(125, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(125, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(125, 16).Value = 2733
$ws.Cells.Item(125, 17).Value = 15
$ws.Cells.Item(125, 18).Value = "Hortaliza"

# Row 126: Zafiro verde, 2021-10-20
$ws.Cells.Item(126, 1).Value = 11
$ws.Cells.Item(126, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(126, 3).Value = "Bíobío"
$ws.Cells.Item(126, 4).Value = 44489
$ws.Cells.Item(126, 5).Value = 8
$ws.Cells.Item(126, 6).Value = 100112002
$ws.Cells.Item(126, 7).Value = "Pimiento"
$ws.Cells.Item(126, 8).Value = "Zafiro verde"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 100
$ws.Cells.Item(126, 11).Value = 33000
$ws.Cells.Item(126, 12).Value = 34000
$ws.Cells.Item(126, 13).Value = 33500
$ws.Cells.Item(126, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(126, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(126, 16).Value = 2233
$ws.Cells.Item(126, 17).Value = 15
$ws.Cells.Item(126, 18).Value = "Hortaliza"
